$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Email" column header (H1) to "Username"
$ws.Range("H1").Value = "Username"

# Move the active selection from H3 to H1
$ws.Range("H1").Select()
